# Fixing integration test Excel regarding hierarchy limit for extensions.
#
# The "Extensions_test" sheet had a RELATION (column D) value on every data
# row (12-26) that exceeded the extension hierarchy limit; those stray
# values are cleared. The sheet also becomes the active/selected tab with
# the selection resting on D13 (where the first removed value used to be),
# while the previously active "ExtensionSchemes_exttest1" sheet loses its
# tabSelected flag and keeps its own prior selection (H2).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Extensions_test")
$ws.Activate()

# Clear the RELATION values that are beyond the allowed extension hierarchy.
$ws.Range("D12:D26").ClearContents()

# Leave the selection on D13, matching the saved workbook state.
$ws.Range("D13").Select()
